# Update "CNKI search terms" sheet: split the old blank spacer row into
# two rows, add two new mini-blocks (CNN / deep learning hanzi recognition)
# mirroring the existing 卷积神经网络/深度学习 rows, each followed by a
# duplicate "title"-row, with a blank separator row between the blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$xlShiftDown   = -4121
$xlContinuous  = 1
$xlThin        = 2
$xlNone        = -4142
$xlEdgeLeft    = 7
$xlEdgeTop     = 8
$xlEdgeBottom  = 9
$xlEdgeRight   = 10
$xlCenter      = -4108
$xlGeneral     = 1

# ---------------------------------------------------------------------
# 1. Split row 18 into two rows: insert a new blank row above the old
#    row 18 (old row 18 -> row 19, old 19..22 -> 20..23, preserving all
#    of their formatting/formulas since it is a true row shift).
# ---------------------------------------------------------------------
$ws.Rows.Item(18).Insert($xlShiftDown)

# Old row 18's formula "=C18+F18" shifted down into B19 verbatim; the
# target sheet has that spacer row with a bare (non-formula) cell, so
# clear it back to empty while keeping its inherited formatting.
$ws.Cells.Item(19, 2).ClearContents()

# ---------------------------------------------------------------------
# 2. New row 18 (the freshly inserted blank row) should look like the
#    other ordinary data rows above it (same borders/number format as
#    row 14, which sits in the same block), just empty.
# ---------------------------------------------------------------------
$ws.Range("A14:J14").Copy()
$ws.Range("A18:J18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
for ($c = 1; $c -le 10; $c++) {
    $ws.Cells.Item(18, $c).ClearContents()
}
$ws.Cells.Item(18, 5).NumberFormat = "General"
$ws.Cells.Item(18, 8).NumberFormat = "General"

# ---------------------------------------------------------------------
# 3. New block 1 (rows 25-26): "CNN 汉字识别" summary row followed by a
#    duplicate of the existing "卷积神经网络 汉字识别" row (row 8).
# ---------------------------------------------------------------------
$ws.Cells.Item(25, 1).Value = "CNN 汉字识别"
$ws.Cells.Item(25, 2).Formula = "=C25+F25"
$ws.Cells.Item(25, 3).Value = 17
$ws.Cells.Item(25, 6).Value = 42

$ws.Cells.Item(26, 1).Value = "卷积神经网络 汉字识别"
$ws.Cells.Item(26, 2).Formula = "=C26+F26"
$ws.Cells.Item(26, 3).Value = 50
$ws.Cells.Item(26, 4).Value = 50
$ws.Cells.Item(26, 5).Formula = '=ROUND(D26/C26*100,2) & "%"'
$ws.Cells.Item(26, 6).Value = 102
$ws.Cells.Item(26, 7).Value = 4
$ws.Cells.Item(26, 8).Formula = '=ROUND(G26/F26*100,2) & "%"'

# ---------------------------------------------------------------------
# 4. Blank separator row 27 (kept empty on purpose).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 5. New block 2 (rows 28-29): "deep learning 汉字识别" summary row
#    followed by a duplicate of the existing "深度学习 汉字识别" row
#    (row 9).
# ---------------------------------------------------------------------
$ws.Cells.Item(28, 1).Value = "deep learning 汉字识别"
$ws.Cells.Item(28, 2).Formula = "=C28+F28"
$ws.Cells.Item(28, 3).Value = 5
$ws.Cells.Item(28, 6).Value = 39

$ws.Cells.Item(29, 1).Value = "深度学习 汉字识别"
$ws.Cells.Item(29, 2).Formula = "=C29+F29"
$ws.Cells.Item(29, 3).Value = 30
$ws.Cells.Item(29, 4).Value = 30
$ws.Cells.Item(29, 5).Formula = '=ROUND(D29/C29*100,2) & "%"'
$ws.Cells.Item(29, 6).Value = 99
$ws.Cells.Item(29, 7).Value = 8
$ws.Cells.Item(29, 8).Formula = '=ROUND(G29/F29*100,2) & "%"'

# ---------------------------------------------------------------------
# 6. Formatting for the two new blocks: mimic the boxed look used
#    elsewhere in the sheet (thin border framing each mini block) and
#    centre the numeric columns, matching rows 19/20 (header+total
#    style) reused here as a "block title" row plus a normal data row.
# ---------------------------------------------------------------------
foreach ($r in 25, 28) {
    for ($c = 1; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
        $cell.Borders.Item($xlEdgeTop).Weight = $xlThin
    }
    $ws.Cells.Item($r, 1).Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
    $ws.Cells.Item($r, 1).Borders.Item($xlEdgeLeft).Weight = $xlThin
    $ws.Cells.Item($r, 3).HorizontalAlignment = $xlCenter
    $ws.Cells.Item($r, 2).HorizontalAlignment = $xlCenter
}

foreach ($r in 26, 29) {
    for ($c = 1; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.HorizontalAlignment = $xlCenter
    }
    $ws.Cells.Item($r, 1).HorizontalAlignment = $xlGeneral
    $ws.Cells.Item($r, 1).Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
    $ws.Cells.Item($r, 1).Borders.Item($xlEdgeLeft).Weight = $xlThin
    $ws.Cells.Item($r, 1).Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $ws.Cells.Item($r, 1).Borders.Item($xlEdgeBottom).Weight = $xlThin
}

$ws.Range("A27:H27").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("A27:H27").Borders.Item($xlEdgeTop).Weight = $xlThin
$ws.Range("A27:H27").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("A27:H27").Borders.Item($xlEdgeBottom).Weight = $xlThin

# ---------------------------------------------------------------------
# 7. Selection marker, matching the author's final cursor position.
# ---------------------------------------------------------------------
$ws.Range("C35").Select()

Write-Output "edit applied"
